$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.891281037216288
$ws.Range("C2").Value = 0.227868235862644
$ws.Range("D2").Value = 0.008829382088016757
$ws.Range("E2").Value = 0.424636301865192
$ws.Range("F2").Value = 0.511900057648603
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("N2").Value = 0.8201990362995204
$ws.Range("O2").Value = 1.626294133104693

$ws.Range("B3").Value = 0.7821342540257206
$ws.Range("C3").Value = 0.2013612742698854
$ws.Range("D3").Value = 0.00796129246657884
$ws.Range("E3").Value = 0.3703224242887444
$ws.Range("F3").Value = 0.4991551764899711
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("N3").Value = 0.8269554030532333
$ws.Range("O3").Value = 1.597674074346656

$ws.Range("B4").Value = 0.715054343910424
$ws.Range("C4").Value = 0.1850036139169049
$ws.Range("D4").Value = 0.007425584329315171
$ws.Range("E4").Value = 0.3370727147557346
$ws.Range("F4").Value = 0.4917844073917266
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("N4").Value = 0.8315100509956181
$ws.Range("O4").Value = 1.581611426959796

$ws.Range("B5").Value = 0.6877035934890046
$ws.Range("C5").Value = 0.1783172446287438
$ws.Range("D5").Value = 0.007206620249466056
$ws.Range("E5").Value = 0.3235461137417133
$ws.Range("F5").Value = 0.4888945045365602
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("N5").Value = 0.8334682874417609
$ws.Range("O5").Value = 1.575443445801824

$ws.Range("B6").Value = 0.683161131083466
$ws.Range("C6").Value = 0.1772057465150567
$ws.Range("D6").Value = 0.007170222232623757
$ws.Range("E6").Value = 0.3213013576355337
$ws.Range("F6").Value = 0.488421491470362
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("N6").Value = 0.8337996237252412
$ws.Range("O6").Value = 1.574442007655904

$ws.Range("B7").Value = 0.7146855423868601
$ws.Range("C7").Value = 0.1849135218971583
$ws.Range("D7").Value = 0.007422633943278356
$ws.Range("E7").Value = 0.3368902001133449
$ws.Range("F7").Value = 0.4917449733671404
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("N7").Value = 0.8315360466860326
$ws.Range("O7").Value = 1.581526716955437

$ws.Range("B8").Value = 0.8536608088396633
$ws.Range("C8").Value = 0.2187458238674935
$ws.Range("D8").Value = 0.008530636135073877
$ws.Range("E8").Value = 0.4058870429058317
$ws.Range("F8").Value = 0.5074109087968353
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("N8").Value = 0.822444341291046
$ws.Range("O8").Value = 1.616111236584089

$ws.Range("B9").Value = 1.125667569260031
$ws.Range("C9").Value = 0.2844324375343206
$ws.Range("D9").Value = 0.01068126729409613
$ws.Range("E9").Value = 0.5420748521896286
$ws.Range("F9").Value = 0.5417675613070116
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("N9").Value = 0.8078376126043167
$ws.Range("O9").Value = 1.696013581328145

$ws.Range("B10").Value = 1.325184521455583
$ws.Range("C10").Value = 0.3322882203825657
$ws.Range("D10").Value = 0.01224697816459752
$ws.Range("E10").Value = 0.642816469174079
$ws.Range("F10").Value = 0.5692687019538027
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("N10").Value = 0.7990697030542435
$ws.Range("O10").Value = 1.762230245765721

$ws.Range("B11").Value = 1.415878935840453
$ws.Range("C11").Value = 0.3539710371430829
$ws.Range("D11").Value = 0.01295597692462991
$ws.Range("E11").Value = 0.6888270794604097
$ws.Range("F11").Value = 0.5822797367671626
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("N11").Value = 0.7955073947389835
$ws.Range("O11").Value = 1.79401746244065

$ws.Range("B12").Value = 1.450212479127003
$ws.Range("C12").Value = 0.3621691174283228
$ws.Range("D12").Value = 0.01322397201867176
$ws.Range("E12").Value = 0.7062788127806243
$ws.Range("F12").Value = 0.5872793665215852
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("N12").Value = 0.7942197524327455
$ws.Range("O12").Value = 1.806296323695108

$ws.Range("B13").Value = 1.442818615038959
$ws.Range("C13").Value = 0.3604040838560536
$ws.Range("D13").Value = 0.01316627646838953
$ws.Range("E13").Value = 0.7025189657318691
$ws.Range("F13").Value = 0.5861993661312965
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("N13").Value = 0.7944943411547314
$ws.Range("O13").Value = 1.803641063067687

$ws.Range("B14").Value = 1.418703791635664
$ws.Range("C14").Value = 0.3546457545478177
$ws.Range("D14").Value = 0.01297803492241911
$ws.Range("E14").Value = 0.6902622597063157
$ws.Range("F14").Value = 0.5826895999395845
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("N14").Value = 0.7954002303527545
$ws.Range("O14").Value = 1.795022793749439

$ws.Range("B15").Value = 1.403931378352922
$ws.Range("C15").Value = 0.3511169470808397
$ws.Range("D15").Value = 0.0128626674887542
$ws.Range("E15").Value = 0.6827584534899103
$ws.Range("F15").Value = 0.5805492449541987
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("N15").Value = 0.7959631012151078
$ws.Range("O15").Value = 1.789775404628529

$ws.Range("B16").Value = 1.319256010063896
$ws.Range("C16").Value = 0.3308694275111179
$ws.Range("D16").Value = 0.01220057630517601
$ws.Range("E16").Value = 0.639813435357766
$ws.Range("F16").Value = 0.5684285290814017
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("N16").Value = 0.7993110881807439
$ws.Range("O16").Value = 1.760186576791938

$ws.Range("B17").Value = 1.267292684093718
$ws.Range("C17").Value = 0.3184257859094828
$ws.Range("D17").Value = 0.01179355792734071
$ws.Range("E17").Value = 0.6135164451768702
$ws.Range("F17").Value = 0.5611215421691469
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("N17").Value = 0.8014741645121291
$ws.Range("O17").Value = 1.74246286157171

$ws.Range("B18").Value = 1.23739851142949
$ws.Range("C18").Value = 0.3112603546215382
$ws.Range("D18").Value = 0.01155914689276472
$ws.Range("E18").Value = 0.5984081774049912
$ws.Range("F18").Value = 0.5569658393738024
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("N18").Value = 0.8027584278971176
$ws.Range("O18").Value = 1.732425220543718

$ws.Range("B19").Value = 1.227275812269681
$ws.Range("C19").Value = 0.3088328648848915
$ws.Range("D19").Value = 0.01147972760001181
$ws.Range("E19").Value = 0.5932956408562973
$ws.Range("F19").Value = 0.5555668567712644
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("N19").Value = 0.8032001470647003
$ws.Range("O19").Value = 1.729053468574108

$ws.Range("B20").Value = 1.272824921058884
$ws.Range("C20").Value = 0.3197512807541898
$ws.Range("D20").Value = 0.01183691742154025
$ws.Range("E20").Value = 0.6163140223165016
$ws.Range("F20").Value = 0.5618945058666469
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("N20").Value = 0.8012397489807483
$ws.Range("O20").Value = 1.744333360078087

$ws.Range("B21").Value = 1.425787190507151
$ws.Range("C21").Value = 0.3563374627162545
$ws.Range("D21").Value = 0.01303333941032747
$ws.Range("E21").Value = 0.6938615595734632
$ws.Range("F21").Value = 0.5837185269568863
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("N21").Value = 0.7951324840674232
$ws.Range("O21").Value = 1.797547609276052

$ws.Range("B22").Value = 1.525695502703115
$ws.Range("C22").Value = 0.3801743374466469
$ws.Range("D22").Value = 0.01381242070084454
$ws.Range("E22").Value = 0.7447106142914066
$ws.Range("F22").Value = 0.598405449573761
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("N22").Value = 0.7914985108803947
$ws.Range("O22").Value = 1.833736282588887

$ws.Range("B23").Value = 1.472378449820724
$ws.Range("C23").Value = 0.3674590195964242
$ws.Range("D23").Value = 0.01339687757999286
$ws.Range("E23").Value = 0.7175554673538045
$ws.Range("F23").Value = 0.5905277838387093
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("N23").Value = 0.7934053094209403
$ws.Range("O23").Value = 1.814291914354243

$ws.Range("B24").Value = 1.270323858226448
$ws.Range("C24").Value = 0.3191520600789204
$ws.Range("D24").Value = 0.01181731587243462
$ws.Range("E24").Value = 0.6150492060403536
$ws.Range("F24").Value = 0.5615449082929018
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("N24").Value = 0.8013456015225415
$ws.Range("O24").Value = 1.743487234541078

$ws.Range("B25").Value = 1.052139364296011
$ws.Range("C25").Value = 0.2667332554221389
$ws.Range("D25").Value = 0.01010192572474011
$ws.Range("E25").Value = 0.5051227919026502
$ws.Range("F25").Value = 0.5320795236704896
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("N25").Value = 0.8114442649019864
$ws.Range("O25").Value = 1.673089269594783
